$d = $word.ActiveDocument
$d.Content.Find.Execute("RFC1 - ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RFC3 - ", 2)
